$d = $word.ActiveDocument

# The target change lives in the "ntpTbl" table style's firstRow conditional
# formatting block (w:tblStylePr w:type="firstRow"): the shading (and the
# borders / vAlign that ship alongside it) applied to the header row need to
# be removed, leaving just the bold/centered run & paragraph formatting.
#
# The Word COM surface exposed here does not implement navigation from a
# Style down to its TableStyle/ConditionalStyle facets (Style.Table /
# Style.Condition always come back empty), so we go around that gap via
# Document.WordOpenXML, which gives us the full package (incl. styles.xml)
# as a single WordprocessingML blob we can edit as text and write back.

$xml = $d.WordOpenXML

$oldBlock = '<w:rsid w:val="00565E19"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:tblPr><w:tblStyleRowBandSize w:val="1"/><w:tblStyleColBandSize w:val="1"/><w:tblInd w:w="0" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblCellMar><w:top w:w="0" w:type="dxa"/><w:left w:w="58" w:type="dxa"/><w:bottom w:w="0" w:type="dxa"/><w:right w:w="58" w:type="dxa"/></w:tblCellMar></w:tblPr><w:trPr><w:cantSplit/></w:trPr><w:tcPr><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:tblStylePr w:type="firstRow"><w:pPr><w:jc w:val="center"/></w:pPr><w:rPr><w:b/><w:bCs/><w:i w:val="0"/><w:iCs w:val="0"/></w:rPr><w:tblPr/><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/><w:vAlign w:val="center"/></w:tcPr></w:tblStylePr>'

$newBlock = '<w:rsid w:val="00B667F1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:tblPr><w:tblStyleRowBandSize w:val="1"/><w:tblStyleColBandSize w:val="1"/><w:tblInd w:w="0" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblCellMar><w:top w:w="0" w:type="dxa"/><w:left w:w="58" w:type="dxa"/><w:bottom w:w="0" w:type="dxa"/><w:right w:w="58" w:type="dxa"/></w:tblCellMar></w:tblPr><w:trPr><w:cantSplit/></w:trPr><w:tcPr><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:tblStylePr w:type="firstRow"><w:pPr><w:jc w:val="center"/></w:pPr><w:rPr><w:b/><w:bCs/><w:i w:val="0"/><w:iCs w:val="0"/></w:rPr></w:tblStylePr>'

if ($xml.Contains($oldBlock)) {
    $xml = $xml.Replace($oldBlock, $newBlock)
    $d.WordOpenXML = $xml
    Write-Host "ntpTbl header-row shading removed."
} else {
    Write-Host "WARNING: expected ntpTbl block not found verbatim; no change applied."
}
